# Add a cumulative "total_runs" column (C) to the "year-runs" sheet,
# mirroring the same cumulative totals already present on Sheet2,
# and update the remembered cell selections on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("year-runs")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Per-year runs (column B, rows 2-17) used to build the running total.
$runs = @(436, 984, 413, 865, 624, 935, 1375, 803, 946, 640, 1095, 606, 805, 747, 771, 1145)

# Header for the new column, formatted like the rest of the normal
# (non-header) data columns on this sheet.
$ws1.Range("C1").Value = "total_runs"
$ws1.Range("C1").NumberFormat = $ws1.Range("B2").NumberFormat
$ws1.Range("C1").Font.Name = $ws1.Range("B2").Font.Name
$ws1.Range("C1").Font.Size = $ws1.Range("B2").Font.Size
$ws1.Range("C1").HorizontalAlignment = $ws1.Range("B2").HorizontalAlignment

$cumulative = 0
for ($i = 0; $i -lt $runs.Length; $i++) {
    $cumulative = $cumulative + $runs[$i]
    $row = $i + 2
    $cell = $ws1.Cells.Item($row, 3)
    $cell.Value = $cumulative
    $cell.NumberFormat = $ws1.Range("B2").NumberFormat
    $cell.Font.Name = $ws1.Range("B2").Font.Name
    $cell.Font.Size = $ws1.Range("B2").Font.Size
    $cell.HorizontalAlignment = $ws1.Range("B2").HorizontalAlignment
}

# Restore the remembered selections recorded in the saved file. Sheet2's
# selection is touched first so that "year-runs" (the tab that was active
# before the edit) ends up selected/active last, leaving the workbook's
# active-sheet state unchanged.
[void]$ws2.Range("B1").Select()
[void]$ws1.Range("D7").Select()
